$d = $word.ActiveDocument

$pairs = @(
    @("88-77=", "88-70="),
    @("40+24=", "85-80="),
    @("46-0=", "13+85="),
    @("9-4=", "69-38="),
    @("64-39=", "67-9="),
    @("90-84=", "41-21="),
    @("13+31=", "29-29="),
    @("37+19=", "28+48="),
    @("70-48=", "9+44="),
    @("9+25=", "0+48="),
    @("95-37=", "66-64="),
    @("11+39=", "3+34="),
    @("95-88=", "27+26="),
    @("10+50=", "38-25="),
    @("21+58=", "2+31="),
    @("29+44=", "13+60="),
    @("78-52=", "38-31="),
    @("39+53=", "94-85="),
    @("77-77=", "29-11="),
    @("34+49=", "90-30="),
    @("45+21=", "57-9="),
    @("22+77=", "56-50="),
    @("33+63=", "97-33="),
    @("93-34=", "45-9="),
    @("32+62=", "47+48="),
    @("70+29=", "93-12="),
    @("19+16=", "50-23="),
    @("3+70=", "51-8="),
    @("66-0=", "38+1="),
    @("32+65=", "76-1="),
    @("11+59=", "74-19="),
    @("18+49=", "19-5="),
    @("63-38=", "30+66="),
    @("38+45=", "96-5="),
    @("29+15=", "23+22="),
    @("6+76=", "39-0="),
    @("58-3=", "82+10="),
    @("79-55=", "29+49="),
    @("43+54=", "11+11="),
    @("40-34=", "1+45="),
    @("99-23=", "79-60="),
    @("68-17=", "10+21="),
    @("72-31=", "23-9="),
    @("21+74=", "3+50="),
    @("74-41=", "76-39="),
    @("4+7=", "41-40="),
    @("35-21=", "64-33="),
    @("7+82=", "28-0="),
    @("82-46=", "8+31="),
    @("73-72=", "96-40="),
    @("66-32=", "16+53="),
    @("22+60=", "62+6="),
    @("31+27=", "57+39="),
    @("45+50=", "10+1="),
    @("20-12=", "39+31="),
    @("85-60=", "95-92="),
    @("44+30=", "86-41="),
    @("27+46=", "61-57="),
    @("48-4=", "71+28="),
    @("4+38=", "6+69="),
    @("32+63=", "56+38="),
    @("19+8=", "46+24="),
    @("39+16=", "26+36="),
    @("61-1=", "35+41="),
    @("42+22=", "34-12="),
    @("57-21=", "10+20="),
    @("30-10=", "88-2="),
    @("67-6=", "26+69="),
    @("23+33=", "82+0="),
    @("48-34=", "51-19="),
    @("72-64=", "14-10="),
    @("6+19=", "1+7="),
    @("45-20=", "91-42="),
    @("7+85=", "56+19="),
    @("76-26=", "38+44="),
    @("85+2=", "13+66="),
    @("79-53=", "28+40="),
    @("74-26=", "13+25="),
    @("81+10=", "28+4="),
    @("92-0=", "62-9="),
    @("35+49=", "56-51="),
    @("92-84=", "24+52="),
    @("45+54=", "27+33="),
    @("3+92=", "37+8="),
    @("92-34=", "92-45="),
    @("20+62=", "81-18="),
    @("90-62=", "25+10="),
    @("45+49=", "79-35="),
    @("42+43=", "30+5="),
    @("32-29=", "89-4="),
    @("33+28=", "84-52="),
    @("21+51=", "1+78="),
    @("9+73=", "67+9="),
    @("50-7=", "80-38="),
    @("10+86=", "27+48="),
    @("95-41=", "48+17="),
    @("8+42=", "42-41="),
    @("24-16=", "83+11="),
    @("77-63=", "4+95="),
    @("2+28=", "13+24=")
)

foreach ($p in $pairs) {
    $old = $p[0]
    $new = $p[1]
    $r = $d.Content
    $found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done."
